$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Restart of the 5-fold experiment: refresh the "train accuracy" column
# (B3:B7) with new run values.
$ws.Range("B3").Value = 0.85046699999999997
$ws.Range("B4").Value = 0.81308400000000003
$ws.Range("B5").Value = 0.82242999999999999
$ws.Range("B5").NumberFormat = "0.000%"
$ws.Range("B6").Value = 0.77570099999999997
$ws.Range("B7").Value = 0.85046699999999997
$ws.Range("B7").NumberFormat = "0.0000%"

# Normalize the 10-fold block's "train accuracy" column formatting/font so
# it matches the rest of the sheet (same font + 0.0000% number format).
$ws.Range("B4").Copy()
$ws.Range("B10:B19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the active selection.
$ws.Range("E10").Select() | Out-Null
